$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.894.52'
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("D3").Value = '4.051.36'
$ws.Range("E3").Value = '  +3.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.79'
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.00'
$ws.Range("E6").Value = '  +2.80%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  +1.51%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.176'
$ws.Range("E10").Value = '  +1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000340'
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.65'
$ws.Range("E12").Value = '  +9.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.72'
$ws.Range("E13").Value = '  +3.60%  '
$ws.Range("D14").Value = '4.702.58'
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("D15").Value = '4.077.60'
$ws.Range("E15").Value = '  +3.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.49'
$ws.Range("E16").Value = '  +8.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.28'
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("E19").Value = '  -1.65%  '
$ws.Range("D20").Value = '71.929.69'
$ws.Range("E20").Value = '  +3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.76'
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.53'
$ws.Range("E22").Value = '  +5.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.77'
$ws.Range("E23").Value = '  +8.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.46'
$ws.Range("E24").Value = '  +7.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.43'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.09'
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.32'
$ws.Range("E27").Value = '  +5.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.24'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.77'
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '704.74'
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  +6.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.52'
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.92'
$ws.Range("E34").Value = '  +15.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '67.34'
$ws.Range("E35").Value = '  -3.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = [string]::Concat('0.0', [char]0x2083, '0906')
$ws.Range("E36").Value = '  +6.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.447'
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.98'
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.61'
$ws.Range("E39").Value = '  +21.42%  '
$ws.Range("E40").Value = '  +4.51%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0488'
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("E44").Value = '  +1.83%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.11'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.54'
$ws.Range("E46").Value = '  +5.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.146'
$ws.Range("E47").Value = '  +3.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.21'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000282'
$ws.Range("E49").Value = '  +23.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.18'
$ws.Range("E50").Value = '  +6.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.34'
$ws.Range("E51").Value = '  +0.31%  '
